$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New notes added to column F for the watering / isotope schedule explanation
$ws.Range("F2").Value = "Watering schedule every Monday & Thursday, unless stated. 100mL each time"
$ws.Range("F3").Value = "Isotope every Wednesday, unless stated. 20mL each time"

# New log rows (serial date numbers, matching the existing A2:A6 cells;
# style is applied below via PasteSpecial so it reuses the existing
# numFmt-14 date style instead of minting a new one)
$ws.Range("A7").Value = 45992
$ws.Range("B7").Value = "did not water"
$ws.Range("C7").Value = "na"
$ws.Range("D7").Value = "worried about overwatering potentially, skipped day"

$ws.Range("A8").Value = 45989
$ws.Range("D8").Value = "some trillium are starting to yellow at tips"

$ws.Range("A9").Value = 45993
$ws.Range("D9").Value = "some trillium showing burn marks"

$ws.Range("B8").Value = "trillium yellowing"
$ws.Range("B9").Value = "trillium burning"

# Apply same date style as existing date cells in column A (reuse the
# existing numFmt-14 style instead of minting a new one)
$ws.Range("A2").Copy()
$ws.Range("A7:A9").PasteSpecial(-4122)

# Widen column B to fit the longer note text
$ws.Columns.Item(2).ColumnWidth = 38.7109375

# Update the active selection like the authored workbook
$ws.Range("B12").Select()
